# Generate Report for Handoff
# Updates status text and timestamps across the Overview / zh-cn / de-de sheets,
# and narrows the "Status" / zh-cn / de-de datetime-adjacent columns.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldOverviewDate = "2016-09-03 03:02:47"
$newOverviewDate = "2016-09-03 03:03:34"

$oldZhDate = "2016-09-03 03:02:43"
$newZhDate = "2016-09-03 03:03:30"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newOverviewDate
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("H2").Value = $newZhDate
$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("H2").Value = $newOverviewDate
$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797
